# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions "Updated cryptos list" commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''61.754.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '''3.390.14'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''579.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '''137.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''3.387.32'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '''7.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '''3.972.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").Value = '''3.382.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''25.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").Value = '''61.850.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '''14.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = '''379.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("D24").Value = '''3.530.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E26").Value = '  +6.72%  '
$ws.Range("D27").Value = '''71.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").Value = '''1.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("D29").Value = '''7.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  +2.38%  '
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = '''3.427.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").Value = '''5.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.07%  '
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").Value = '''6.87'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("D40").Value = '''165.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("E42").Value = '  +2.85%  '
$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").Value = '''1.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("E45").Value = '  +7.07%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''25.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.00%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '''4.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").Value = '''41.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("D50").Value = '''22.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("D51").Value = '''2.349.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.37%  '
